$wb = $excel.ActiveWorkbook

# --- NormalSearch: update the selected range (no change to data) ---
$normalSearch = $wb.Worksheets.Item("NormalSearch")
$normalSearch.Range("A1:A6").Select()

# --- Insert a new sheet "NormalSearch(2)" right after "NormalSearch" ---
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $normalSearch)
$newSheet.Name = "NormalSearch(2)"
$newSheet.Range("A1").Value = "searchkeyword"
$newSheet.Range("A1").Font.Bold = $true
$newSheet.Range("A2").Value = "companies"
$newSheet.Range("A3").Value = "televisions"
$newSheet.Range("A4").Value = "ninjas"
$newSheet.Range("D25").Select()

# --- FeelingLuckySearch: update row 2 to hold the new search term/value,
#     drop row 3, and resize column A to fit the new content ---
$feelingLucky = $wb.Worksheets.Item("FeelingLuckySearch")
$feelingLucky.Range("A2").Value = "searchkeyword"
$feelingLucky.Range("B2").Value = "im feeling lucky"
$feelingLucky.Range("A3:B3").ClearContents()
$feelingLucky.Columns("A:A").ColumnWidth = 13.71

# --- Make FeelingLuckySearch the active sheet/cell, as it is now the one being configured ---
$feelingLucky.Activate()
$feelingLucky.Range("B3").Select()
